# Add a new "test_case_qualifiers" column to the TestCase-family sheets,
# inserted immediately before the existing "test_case_input_id" column.
#
# Affected sheets and the column currently holding "test_case_input_id":
#   TestCase                          -> K
#   AcceptanceTestCase                -> K
#   QuantitativeTestCase              -> K
#   ComplianceTestCase                -> M
#   KnowledgeGraphNavigationTestCase  -> K
#   OneHopTestCase                    -> K

$wb = $excel.ActiveWorkbook

$sheetsToColumn = @{
    "TestCase"                         = "K"
    "AcceptanceTestCase"                = "K"
    "QuantitativeTestCase"              = "K"
    "ComplianceTestCase"                = "M"
    "KnowledgeGraphNavigationTestCase"  = "K"
    "OneHopTestCase"                    = "K"
}

foreach ($sheetName in $sheetsToColumn.Keys) {
    $col = $sheetsToColumn[$sheetName]
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new blank column before the "test_case_input_id" column,
    # shifting it (and everything after it) one column to the right.
    $ws.Range($col + "1").EntireColumn.Insert()

    # Populate the header of the newly inserted column.
    $ws.Range($col + "1").Value = "test_case_qualifiers"
}
